# Update cryptocurrency price/volume data as scraped on Tue Oct 31 02:42:20 UTC 2023
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "34.444.66"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.22%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.810.03"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.17%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "228.27"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.61%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.579"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.12%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "35.77"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +8.41%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.303"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.65%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0696"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.96%  "
$ws.Range("E11").Value = "  +1.96%  "
$ws.Range("E12").Value = "  +1.26%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.50"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.53%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.809.63"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.83%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.647"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.76%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.54"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +5.59%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "34.469.68"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.31%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "69.33"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.98%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "247.34"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.74%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0800"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.18%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.55"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.12%  "
$ws.Range("E22").Value = "  +0.01%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.22"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.34%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "171.72"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.89%  "
$ws.Range("E25").Value = "  +3.20%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.97"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +8.36%  "
$ws.Range("E27").Value = "  +2.81%  "
$ws.Range("E28").Value = "  +3.82%  "
$ws.Range("E29").Value = "  +0.07%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.08"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.01%  "
$ws.Range("B31").Value = "Hedera"
$ws.Range("C31").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0535"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.58%  "
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.88"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.17%  "
$ws.Range("E33").Value = "  +1.06%  "
$ws.Range("E34").Value = "  +1.47%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.399.81"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.02%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.679"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.06%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.48"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.46%  "
$ws.Range("E38").Value = "  +0.42%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0191"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.58%  "
$ws.Range("B40").Value = "WEMIXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.26"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +14.02%  "
$ws.Range("B41").Value = "Aave"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "83.30"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.38%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.968"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.64%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.83"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.83%  "
$ws.Range("E44").Value = "  +0.09%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.40"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.41%  "
$ws.Range("E46").Value = "  -1.35%  "
$ws.Range("E47").Value = "  -4.18%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.973.07"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.10%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "105.83"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.53%  "
$ws.Range("E50").Value = "  +0.07%  "
$ws.Range("E51").Value = "  +1.15%  "
